$d = $word.ActiveDocument

# 1. Title: "Disadvantages of e-coaching" -> "Disadvantages of eCoaching"
$d.Content.Find.Execute("Disadvantages of e-coaching", $true, $false, $false, $false, $false, $true, 1, $false, "Disadvantages of eCoaching", 2) | Out-Null

# 2. "E-coaching contains a lot of advantages" -> "eCoaching contains a lot of advantages"
$d.Content.Find.Execute("E-coaching contains a lot of advantages", $true, $false, $false, $false, $false, $true, 1, $false, "eCoaching contains a lot of advantages", 2) | Out-Null

# 3. "will sooner disengage from e-coaches." -> "will sooner disengage from eCoaches."
$d.Content.Find.Execute("will sooner disengage from e-coaches.", $true, $false, $false, $false, $false, $true, 1, $false, "will sooner disengage from eCoaches.", 2) | Out-Null

# 4. "A long with the commitment issue, e-coaches are also harder" -> "...eCoaches are also harder"
$d.Content.Find.Execute("A long with the commitment issue, e-coaches are also harder", $true, $false, $false, $false, $false, $true, 1, $false, "A long with the commitment issue, eCoaches are also harder", 2) | Out-Null

# 5. "use the e-coaches fully as they were intended. It is possible" -> "use the eCoach fully as they were intended to use. It is possible"
$d.Content.Find.Execute("use the e-coaches fully as they were intended. It is possible", $true, $false, $false, $false, $false, $true, 1, $false, "use the eCoach fully as they were intended to use. It is possible", 2) | Out-Null

# 6. "do not know how to work with. Interacting" -> "do not know how to work with the program. Interacting"
$d.Content.Find.Execute("do not know how to work with. Interacting", $true, $false, $false, $false, $false, $true, 1, $false, "do not know how to work with the program. Interacting", 2) | Out-Null

# 7. "will slow down progress to get over" -> "will slow down the progress to get over"
$d.Content.Find.Execute("will slow down progress to get over", $true, $false, $false, $false, $false, $true, 1, $false, "will slow down the progress to get over", 2) | Out-Null

# 8. "E-coaches are not as good as people yet." -> "eCoaches are not as good as people, yet."
$d.Content.Find.Execute("E-coaches are not as good as people yet.", $true, $false, $false, $false, $false, $true, 1, $false, "eCoaches are not as good as people, yet.", 2) | Out-Null

# 9. Insert the _GoBack bookmark right before " yet." (i.e. right after the comma we just added)
$full = $d.Content.Text
$pos = $full.IndexOf("eCoaches are not as good as people, yet.")
$bmPos = $pos + "eCoaches are not as good as people,".Length
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# 10. Citation: remove stale proofErr spell-check runs and merge into a single run.
#     A same-text replace is a no-op in this engine, so round-trip through a
#     placeholder to force the run to be rebuilt.
$full = $d.Content.Text
$pos = $full.IndexOf("David Clutterbuck")
$len = "David Clutterbuck & Zulfi Hussain. 2010.".Length
$r = $d.Range($pos, $pos + $len)
$r.Text = "IRONNATIVE_TMP_PLACEHOLDER"
$full2 = $d.Content.Text
$pos2 = $full2.IndexOf("IRONNATIVE_TMP_PLACEHOLDER")
$r2 = $d.Range($pos2, $pos2 + "IRONNATIVE_TMP_PLACEHOLDER".Length)
$r2.Text = "David Clutterbuck & Zulfi Hussain. 2010."
